$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '29.212.80'
$ws.Range("E2").Value = '  -0.04%  '
$ws.Range("D3").Value = '1.849.50'
$ws.Range("E3").Value = '  -0.27%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.9995'
$ws.Range("E4").Value = '  -0.11%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '245.78'
$ws.Range("E5").Value = '  +1.93%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.6990'
$ws.Range("E6").Value = '  -0.03%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.9999'
$ws.Range("E7").Value = '  -0.11%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.07723'
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.3066'
$ws.Range("E9").Value = '  -0.77%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '23.52'
$ws.Range("E10").Value = '  -0.92%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.07816'
$ws.Range("E11").Value = '  +0.22%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '92.90'
$ws.Range("E12").Value = '  +0.96%  '
$ws.Range("D13").Value = '1.847.35'
$ws.Range("E13").Value = '  -0.69%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '5.128'
$ws.Range("E14").Value = '  +0.72%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.6863'
$ws.Range("E15").Value = '  -0.06%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '6.627'
$ws.Range("E16").Value = '  +2.35%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.000008321'
$ws.Range("E17").Value = '  -0.81%  '
$ws.Range("D18").Value = '29.206.35'
$ws.Range("E18").Value = '  -0.16%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '241.62'
$ws.Range("E19").Value = '  -3.04%  '
$ws.Range("D20").Value = '2.088.29'
$ws.Range("E20").Value = '  -1.79%  '
$ws.Range("E21").Value = '  -0.74%  '
$ws.Range("E22").Value = '  -0.08%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '7.528'
$ws.Range("E23").Value = '  +0.24%  '
$ws.Range("E24").Value = '  -0.09%  '
$ws.Range("E25").Value = '  -1.16%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '159.03'
$ws.Range("E26").Value = '  -0.78%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '8.828'
$ws.Range("E27").Value = '  -0.06%  '
$ws.Range("E28").Value = '  -1.12%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '1.542'
$ws.Range("E29").Value = '  -1.11%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '4.230'
$ws.Range("E30").Value = '  +0.06%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '4.183'
$ws.Range("E31").Value = '  -0.32%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '1.197'
$ws.Range("E32").Value = '  +0.43%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.05123'
$ws.Range("E33").Value = '  -1.22%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.7932'
$ws.Range("E34").Value = '  +4.58%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.901'
$ws.Range("E35").Value = '  +3.27%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '1.147'
$ws.Range("E36").Value = '  -1.05%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '2.694'
$ws.Range("E37").Value = '  -0.59%  '
$ws.Range("D38").Value = '1.323.22'
$ws.Range("E38").Value = '  +7.83%  '
$ws.Range("E39").Value = '  +0.87%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '2.717'
$ws.Range("E40").Value = '  -0.30%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.9497'
$ws.Range("E41").Value = '  +6.07%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '6.014'
$ws.Range("E42").Value = '  +7.92%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '107.33'
$ws.Range("E43").Value = '  -2.23%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '1.000'
$ws.Range("E44").Value = '  -0.04%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '9.726'
$ws.Range("E45").Value = '  +2.59%  '
$ws.Range("D46").Value = '1.989.76'
$ws.Range("E46").Value = '  -1.18%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.5184'
$ws.Range("E47").Value = '  +0.05%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '64.17'
$ws.Range("E48").Value = '  -1.87%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.766'
$ws.Range("E49").Value = '  +1.31%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.00000000118'
$ws.Range("E50").Value = '  -2.84%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '7.005'
$ws.Range("E51").Value = '  +0.20%  '
